$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 100, 102 and 103 got their match-data (everything except the
#    running "id" counter in column A, and the Div/Div Original Name/Date
#    columns C:E which were already identical across the three rows) cycled:
#      new row100 <- old row102
#      new row102 <- old row103
#      new row103 <- old row100
# ---------------------------------------------------------------------------

# --- Row 100 (<- old row 102 values) ---
$ws.Range("B100").Value = 6732834
$ws.Range("F100").Value = "Panevezys"
$ws.Range("G100").Value = "FK Dziugas Telsiai"
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = "D"
$ws.Range("K100").Value = 1.25
$ws.Range("L100").Value = 5.5
$ws.Range("M100").Value = 7.5
$ws.Range("N100").Value = 1.45
$ws.Range("O100").Value = 4.5
$ws.Range("P100").Value = 5
$ws.Range("Q100").Value = -1
$ws.Range("R100").Value = 1.775
$ws.Range("S100").Value = 2.025
$ws.Range("T100").Value = 2.5
$ws.Range("U100").Value = 1.875
$ws.Range("V100").Value = 1.925
$ws.Range("W100").Value = -1
$ws.Range("X100").Value = 3.5
$ws.Range("Y100").Value = -1
$ws.Range("Z100").Value = -1
$ws.Range("AA100").Value = 1.025
$ws.Range("AB100").Value = -1
$ws.Range("AC100").Value = 0.925

# --- Row 102 (<- old row 103 values) ---
$ws.Range("B102").Value = 6732727
$ws.Range("F102").Value = "FK Zalgiris Vilnius"
$ws.Range("G102").Value = "FK Dainava Alytus"
$ws.Range("H102").Value = 1
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = "H"
$ws.Range("K102").Value = 1.285
$ws.Range("L102").Value = 5.5
$ws.Range("M102").Value = 6.5
$ws.Range("N102").Value = 1.3
$ws.Range("O102").Value = 5.5
$ws.Range("P102").Value = 6
$ws.Range("Q102").Value = -1.5
$ws.Range("R102").Value = 1.9
$ws.Range("S102").Value = 1.9
$ws.Range("T102").Value = 2.75
$ws.Range("U102").Value = 1.8
$ws.Range("V102").Value = 2
$ws.Range("W102").Value = 0.3
$ws.Range("X102").Value = -1
$ws.Range("Y102").Value = -1
$ws.Range("Z102").Value = -1
$ws.Range("AA102").Value = 0.8999999999999999
$ws.Range("AB102").Value = -1
$ws.Range("AC102").Value = 1

# --- Row 103 (<- old row 100 values) ---
$ws.Range("B103").Value = 7465686
$ws.Range("F103").Value = "FK Kauno Zalgiris"
$ws.Range("G103").Value = "Hegelmann Litauen"
$ws.Range("H103").Value = 4
$ws.Range("I103").Value = 2
$ws.Range("J103").Value = "H"
$ws.Range("K103").Value = 2.3
$ws.Range("L103").Value = 4
$ws.Range("M103").Value = 2.3
$ws.Range("N103").Value = 2.55
$ws.Range("O103").Value = 4
$ws.Range("P103").Value = 2.2
$ws.Range("Q103").Value = 0.25
$ws.Range("R103").Value = 1.8
$ws.Range("S103").Value = 2
$ws.Range("T103").Value = 2.75
$ws.Range("U103").Value = 1.85
$ws.Range("V103").Value = 1.95
$ws.Range("W103").Value = 1.55
$ws.Range("X103").Value = -1
$ws.Range("Y103").Value = -1
$ws.Range("Z103").Value = 0.8
$ws.Range("AA103").Value = -1
$ws.Range("AB103").Value = 0.8500000000000001
$ws.Range("AC103").Value = -1

# ---------------------------------------------------------------------------
# 2) Append two new fixtures as rows 135 and 136.
# ---------------------------------------------------------------------------

# --- Row 135 : FK Dziugas Telsiai related game, already played ---
$ws.Range("A135").Value = 133
$ws.Range("B135").Value = 7862921
$ws.Range("C135").Value = "Lithuania A Lyga"
$ws.Range("D135").Value = "Lithuania A Lyga"
$ws.Range("E135").Value = 45391.5
$ws.Range("F135").Value = "Hegelmann Litauen"
$ws.Range("G135").Value = "FK Dainava Alytus"
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 1
$ws.Range("J135").Value = "A"
$ws.Range("K135").Value = 1.75
$ws.Range("L135").Value = 3.25
$ws.Range("M135").Value = 4.2
$ws.Range("N135").Value = 1.45
$ws.Range("O135").Value = 3.6
$ws.Range("P135").Value = 6
$ws.Range("Q135").Value = -1
$ws.Range("R135").Value = 1.875
$ws.Range("S135").Value = 1.925
$ws.Range("T135").Value = 2.5
$ws.Range("U135").Value = 1.9
$ws.Range("V135").Value = 1.9
$ws.Range("W135").Value = -1
$ws.Range("X135").Value = -1
$ws.Range("Y135").Value = 5
$ws.Range("Z135").Value = -1
$ws.Range("AA135").Value = 0.925
$ws.Range("AB135").Value = -1
$ws.Range("AC135").Value = 0.8999999999999999

# --- Row 136 : upcoming fixture, result/PL columns not available yet ---
$ws.Range("A136").Value = 134
$ws.Range("B136").Value = 7862922
$ws.Range("C136").Value = "Lithuania A Lyga"
$ws.Range("D136").Value = "Lithuania A Lyga"
$ws.Range("E136").Value = 45392.5
$ws.Range("F136").Value = "FK Siauliai"
$ws.Range("G136").Value = "Panevezys"
$ws.Range("K136").Value = 2.7
$ws.Range("L136").Value = 3
$ws.Range("M136").Value = 2.5
$ws.Range("N136").Value = 2.6
$ws.Range("O136").Value = 3
$ws.Range("P136").Value = 2.6
$ws.Range("Q136").Value = 0
$ws.Range("R136").Value = 1.9
$ws.Range("S136").Value = 1.9
$ws.Range("T136").Value = 2
$ws.Range("U136").Value = 1.9
$ws.Range("V136").Value = 1.9
$ws.Range("W136").Value = 0
$ws.Range("X136").Value = 0
$ws.Range("Y136").Value = 0
$ws.Range("Z136").Value = 0
$ws.Range("AA136").Value = 0

# Column A (bold/centered/bordered) and column E (date format) carry special
# styling in this sheet; copy it from the previous data row (134) onto the
# two freshly appended rows.
$ws.Range("A134").Copy() | Out-Null
$ws.Range("A135:A136").PasteSpecial(-4122) | Out-Null
$ws.Range("E134").Copy() | Out-Null
$ws.Range("E135:E136").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
